$wb = $excel.ActiveWorkbook

# --- Before state --------------------------------------------------------
# Worksheets(1) = "hotel_info"  -> header row A1:I1, one data row A2:I2
# Worksheets(2) = "review_info" -> header row A1:Y1 (25 cols), no data row
#
# --- After state (per diff) ----------------------------------------------
# Physical sheet #1 now carries the review_info content (25 header columns,
#   still no data row) and is renamed "review_info".
# Physical sheet #2 now carries the hotel_info content (header + 1 data
#   row) with a new "State" column inserted right after "Hotel_Name"
#   (value "Louisiana" in the data row) and is renamed "hotel_info".
# Tab order keeps review_info before hotel_info (same physical order as
# before - only names + contents are swapped).
#
# NB: `Range.Copy(destRange)` is used (instead of re-typing `.Value()`
# results) so numeric-looking text such as "559"/"71"/"572" keeps its
# original shared-string/text cell type instead of being re-inferred as a
# number by a plain PowerShell value assignment.
# ---------------------------------------------------------------------------

$s1 = $wb.Worksheets.Item(1)   # currently "hotel_info"
$s2 = $wb.Worksheets.Item(2)   # currently "review_info"

# --- Stage everything we still need onto scratch rows on the OTHER sheet
# before anything gets cleared/overwritten -> both stashes read pristine,
# untouched source data. ----------------------------------------------------

# Stash review_info's own header row (from s2) onto s1's scratch row 100.
$s2.Range("A1:Y1").Copy($s1.Range("A100:Y100"))

# Stash hotel_info's own header + data rows (from s1) onto s2's scratch
# rows 100/101.
$s1.Range("A1:I1").Copy($s2.Range("A100:I100"))
$s1.Range("A2:I2").Copy($s2.Range("A101:I101"))

# --- Clear the original content (scratch rows are untouched) -------------
$s1.Range("A1:I2").Clear()
$s2.Range("A1:Y1").Clear()

# --- Write final content of the sheet that becomes "review_info" ---------
$s1.Range("A100:Y100").Copy($s1.Range("A1:Y1"))

# --- Write final content of the sheet that becomes "hotel_info" ----------
# Column map: src A..I (STR,Hotel_Name,City,Zip,TA_ReviewURL,
#   Tripadvisor_Hotel_Name,English_Reviews_num,Local_Rank,Total_Reviews_num)
#   -> dst A,B,D,E,F,G,H,I,J  (C is the new "State" column)
$srcCols = @(1,2,3,4,5,6,7,8,9)
$dstCols = @(1,2,4,5,6,7,8,9,10)
for ($i = 0; $i -lt $srcCols.Length; $i++) {
    $sc = $srcCols[$i]
    $dc = $dstCols[$i]
    $s2.Cells.Item(100, $sc).Copy($s2.Cells.Item(1, $dc))
    $s2.Cells.Item(101, $sc).Copy($s2.Cells.Item(2, $dc))
}
$s2.Cells.Item(1, 3).Value = "State"
$s2.Cells.Item(2, 3).Value = "Louisiana"

# --- Remove scratch rows ---------------------------------------------------
$s1.Range("A100:Y100").Clear()
$s2.Range("A100:I101").Clear()

# --- Rename (temp name avoids a transient name collision) ----------------
$s1.Name = "review_info_tmp"
$s2.Name = "hotel_info"
$s1.Name = "review_info"
